$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Efna5"
$ws.Range("C2").Value = "Ephb6"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.030023666666667
$ws.Range("H2").Value = 6.090071
$ws.Range("I2").Value = 0.8776223887075381
$ws.Range("J2").Value = 0.8776223887075382
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.167255
$ws.Range("N2").Value = 0.501765
$ws.Range("O2").Value = 0.1064977723687984
$ws.Range("P2").Value = 0.1064977723687984
$ws.Range("Q2").Value = 0.3395316083683333
$ws.Range("R2").Value = 3.055784475315
$ws.Range("S2").Value = 0.09346482937833647
$ws.Range("T2").Value = 0.09346482937833649

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Efna5"
$ws.Range("C3").Value = "Ephb6"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.030023666666667
$ws.Range("H3").Value = 6.090071
$ws.Range("I3").Value = 0.8776223887075381
$ws.Range("J3").Value = 0.8776223887075382
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.429376
$ws.Range("N3").Value = 1.288128
$ws.Range("O3").Value = 0.2734004215636314
$ws.Range("P3").Value = 0.2734004215636314
$ws.Range("Q3").Value = 0.8716434418986666
$ws.Range("R3").Value = 7.844790977088
$ws.Range("S3").Value = 0.2399423310463221
$ws.Range("T3").Value = 0.2399423310463221

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Efna5"
$ws.Range("C4").Value = "Ephb6"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.030023666666667
$ws.Range("H4").Value = 6.090071
$ws.Range("I4").Value = 0.8776223887075381
$ws.Range("J4").Value = 0.8776223887075382
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02801
$ws.Range("N4").Value = 0.08402999999999999
$ws.Range("O4").Value = 0.01783505787001908
$ws.Range("P4").Value = 0.01783505787001909
$ws.Range("Q4").Value = 0.05686096290333333
$ws.Range("R4").Value = 0.5117486661299999
$ws.Range("S4").Value = 0.01565244609062332
$ws.Range("T4").Value = 0.01565244609062333

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Efna5"
$ws.Range("C5").Value = "Ephb6"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.030023666666667
$ws.Range("H5").Value = 6.090071
$ws.Range("I5").Value = 0.8776223887075381
$ws.Range("J5").Value = 0.8776223887075382
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.9458613333333332
$ws.Range("N5").Value = 2.837584
$ws.Range("O5").Value = 0.6022667481975512
$ws.Range("P5").Value = 0.6022667481975512
$ws.Range("Q5").Value = 1.920120892051555
$ws.Range("R5").Value = 17.281088028464
$ws.Range("S5").Value = 0.5285627821922563
$ws.Range("T5").Value = 0.5285627821922563

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Efna5"
$ws.Range("C6").Value = "Ephb6"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.283071
$ws.Range("H6").Value = 0.849213
$ws.Range("I6").Value = 0.1223776112924619
$ws.Range("J6").Value = 0.1223776112924619
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.167255
$ws.Range("N6").Value = 0.501765
$ws.Range("O6").Value = 0.1064977723687984
$ws.Range("P6").Value = 0.1064977723687984
$ws.Range("Q6").Value = 0.04734504010500001
$ws.Range("R6").Value = 0.426105360945
$ws.Range("S6").Value = 0.01303294299046189
$ws.Range("T6").Value = 0.01303294299046189

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Efna5"
$ws.Range("C7").Value = "Ephb6"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.283071
$ws.Range("H7").Value = 0.849213
$ws.Range("I7").Value = 0.1223776112924619
$ws.Range("J7").Value = 0.1223776112924619
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.429376
$ws.Range("N7").Value = 1.288128
$ws.Range("O7").Value = 0.2734004215636314
$ws.Range("P7").Value = 0.2734004215636314
$ws.Range("Q7").Value = 0.121543893696
$ws.Range("R7").Value = 1.093895043264
$ws.Range("S7").Value = 0.03345809051730929
$ws.Range("T7").Value = 0.03345809051730929

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Efna5"
$ws.Range("C8").Value = "Ephb6"
$ws.Range("D8").Value = "M1"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.283071
$ws.Range("H8").Value = 0.849213
$ws.Range("I8").Value = 0.1223776112924619
$ws.Range("J8").Value = 0.1223776112924619
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.02801
$ws.Range("N8").Value = 0.08402999999999999
$ws.Range("O8").Value = 0.01783505787001908
$ws.Range("P8").Value = 0.01783505787001909
$ws.Range("Q8").Value = 0.007928818709999999
$ws.Range("R8").Value = 0.07135936839
$ws.Range("S8").Value = 0.002182611779395758
$ws.Range("T8").Value = 0.002182611779395759

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Efna5"
$ws.Range("C9").Value = "Ephb6"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.283071
$ws.Range("H9").Value = 0.849213
$ws.Range("I9").Value = 0.1223776112924619
$ws.Range("J9").Value = 0.1223776112924619
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.9458613333333332
$ws.Range("N9").Value = 2.837584
$ws.Range("O9").Value = 0.6022667481975512
$ws.Range("P9").Value = 0.6022667481975512
$ws.Range("Q9").Value = 0.267745913488
$ws.Range("R9").Value = 2.409713221392
$ws.Range("S9").Value = 0.07370396600529494
$ws.Range("T9").Value = 0.07370396600529494

